$d = $word.ActiveDocument

# Locate the paragraph that still carries the old, unfinished
# "@DatabaseFactory, @DatabaseMySql, AccountDAO" text (split across two
# runs in the original document).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*@DatabaseFactory, @DatabaseMySql*") {
        $target = $p
    }
}

# Append a brand-new (empty) paragraph right after it - Word clones the
# paragraph-mark formatting automatically, so it inherits the same
# Times New Roman / 26-half-point run formatting.
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$newPara = $target.Next()

# That new paragraph becomes the relocated "@DatabaseFactory, @DatabaseMySql" line.
$newPara.Range.InsertBefore("@DatabaseFactory, @DatabaseMySql")

# The original paragraph turns into the new "Folder Dao" label - replace
# its whole (pre-paragraph-mark) text, which drops the trailing
# ", AccountDAO" run along with it.
$pr = $target.Range
[void]$pr.MoveEnd(1, -1)
$pr.Text = "Folder Dao"
